$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "name" column (C) with the person associated with each
# existing email/password row. Order of assignment matches the order new
# shared strings were introduced in the source workbook.
$ws.Range("C4").Value = "Florian"
$ws.Range("C1").Value = "SODKI"
$ws.Range("C2").Value = "Viktoriia"
$ws.Range("C3").Value = "Haroun"
$ws.Range("C7").Value = "Salma"
$ws.Range("C5").Value = "Said"
$ws.Range("C6").Value = "Imed"
$ws.Range("C8").Value = "Yacine"

# Remove the very last account (row 9: e3u3@utopios.solutions / VCj*2Eqj),
# including its hyperlink. Hyperlinks.Delete() on this runtime clears every
# hyperlink on the sheet, so we drop them all and re-create the ones that
# must survive.
$ws.Hyperlinks.Delete()
$ws.Range("A9").Value = $null
$ws.Range("B9").Value = $null

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:e2u2@utopios.solutions")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:e1u2@utopios.solutions")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:e2u1@utopios.solutions")
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:e2u3@utopios.solutions")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:e3u1@utopios.solutions")
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:e3u2@utopios.solutions")

# Restore the hyperlink cell style (Hyperlinks.Add() bumps the style table,
# so re-apply the named "Lien hypertexte" style explicitly to keep using
# the original shared cell format).
$ws.Range("A2").Style = "Lien hypertexte"
$ws.Range("A3").Style = "Lien hypertexte"
$ws.Range("A5").Style = "Lien hypertexte"
$ws.Range("A6").Style = "Lien hypertexte"
$ws.Range("A7").Style = "Lien hypertexte"
$ws.Range("A8").Style = "Lien hypertexte"
$ws.Range("A9").Style = "Lien hypertexte"

# Adjust column widths to fit the new layout. (The values below are chosen
# so that, after this runtime's internal char->pixel width quantization,
# the saved column widths land on 42.83203125 / 30.5 / 21 -- matching what
# Excel itself persists in xl/worksheets/sheet1.xml.)
$ws.Columns.Item(1).ColumnWidth = 42.0
$ws.Columns.Item(2).ColumnWidth = 29.666666666666668
$ws.Columns.Item(3).ColumnWidth = 20.166666666666668

# Update the selected cell to reflect the new active cell position
$ws.Range("C9").Select()
